$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three obsolete rows (old rows 5-7, MuSCs as sender) entirely
$ws.Rows("5:7").Delete()

# Clear the existing data rows (2-4) so stale shared-string references are dropped
$ws.Range("A2:T4").Clear()

# Re-populate rows 2-4 with the refreshed TPM-based values
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Il1b"
$ws.Cells.Item(2, 3).Value = "Il1rap"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.1046376666666667
$ws.Cells.Item(2, 8).Value = 0.313913
$ws.Cells.Item(2, 9).Value = 1
$ws.Cells.Item(2, 10).Value = 1
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 1.135923333333333
$ws.Cells.Item(2, 14).Value = 3.40777
$ws.Cells.Item(2, 15).Value = 0.06998805686568385
$ws.Cells.Item(2, 16).Value = 0.06998805686568385
$ws.Cells.Item(2, 17).Value = 0.1188603671122222
$ws.Cells.Item(2, 18).Value = 1.06974330401
$ws.Cells.Item(2, 19).Value = 0.06998805686568385
$ws.Cells.Item(2, 20).Value = 0.06998805686568385
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Il1b"
$ws.Cells.Item(3, 3).Value = "Il1rap"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.1046376666666667
$ws.Cells.Item(3, 8).Value = 0.313913
$ws.Cells.Item(3, 9).Value = 1
$ws.Cells.Item(3, 10).Value = 1
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 7.095305
$ws.Cells.Item(3, 14).Value = 21.285915
$ws.Cells.Item(3, 15).Value = 0.4371656037403091
$ws.Cells.Item(3, 16).Value = 0.437165603740309
$ws.Cells.Item(3, 17).Value = 0.7424361594883333
$ws.Cells.Item(3, 18).Value = 6.681925435395
$ws.Cells.Item(3, 19).Value = 0.4371656037403091
$ws.Cells.Item(3, 20).Value = 0.437165603740309
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Il1b"
$ws.Cells.Item(4, 3).Value = "Il1rap"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.1046376666666667
$ws.Cells.Item(4, 8).Value = 0.313913
$ws.Cells.Item(4, 9).Value = 1
$ws.Cells.Item(4, 10).Value = 1
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 7.999016999999999
$ws.Cells.Item(4, 14).Value = 23.997051
$ws.Cells.Item(4, 15).Value = 0.4928463393940071
$ws.Cells.Item(4, 16).Value = 0.4928463393940071
$ws.Cells.Item(4, 17).Value = 0.8369984745069999
$ws.Cells.Item(4, 18).Value = 7.532986270563
$ws.Cells.Item(4, 19).Value = 0.4928463393940071
$ws.Cells.Item(4, 20).Value = 0.4928463393940071
